$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 1780.25
$ws.Range("J100").Value = 1968.3334
$ws.Range("L100").Value = 1968.3334
$ws.Range("N100").Value = -3050.3334
# Row 112
$ws.Range("H112").Value = 980.375
$ws.Range("J112").Value = 980.375
$ws.Range("L112").Value = 2941.125
$ws.Range("N112").Value = -5157.125
# Row 116
$ws.Range("H116").Value = 3444
$ws.Range("I116").Value = 3444
$ws.Range("K116").Value = 3444
$ws.Range("M116").Value = -2
# Row 131
$ws.Range("H131").Value = 9814.333000000001
$ws.Range("I131").Value = 1451.75
$ws.Range("J131").Value = 14960.538
$ws.Range("K131").Value = 4355.25
$ws.Range("L131").Value = 44881.614
$ws.Range("M131").Value = 684.75
$ws.Range("N131").Value = -54961.614
# Row 138
$ws.Range("H138").Value = 3835.7158
$ws.Range("I138").Value = 3124.4092
$ws.Range("J138").Value = 4072.818
$ws.Range("K138").Value = 9373.2276
$ws.Range("L138").Value = 12218.454
$ws.Range("M138").Value = -4233.2276
$ws.Range("N138").Value = -22498.454

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2121.6875
$ws.Range("I2").Value = 994.8
$ws.Range("J2").Value = 3999.8333
$ws.Range("K2").Value = 994.8
$ws.Range("L2").Value = 3999.8333
$ws.Range("M2").Value = -881.8
$ws.Range("N2").Value = -4225.8333
# Row 17
$ws.Range("H17").Value = 20165.166
$ws.Range("J17").Value = 20333.666
$ws.Range("L17").Value = 20333.666
$ws.Range("N17").Value = -20679.666
# Row 45
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1623
# Row 61
$ws.Range("H61").Value = 4226.1
$ws.Range("I61").Value = 1758.5
$ws.Range("K61").Value = 1758.5
$ws.Range("M61").Value = -1546.5
# Row 110
$ws.Range("H110").Value = 1308.4286
$ws.Range("I110").Value = 1308.4286
$ws.Range("K110").Value = 1308.4286
$ws.Range("M110").Value = 736.5714
# Row 116
$ws.Range("H116").Value = 2121.6875
$ws.Range("I116").Value = 994.8
$ws.Range("J116").Value = 3999.8333
$ws.Range("K116").Value = 994.8
$ws.Range("L116").Value = 3999.8333
$ws.Range("M116").Value = 1299.2
$ws.Range("N116").Value = -8587.8333
# Row 136
$ws.Range("H136").Value = 4226.1
$ws.Range("I136").Value = 1758.5
$ws.Range("K136").Value = 5275.5
$ws.Range("M136").Value = -2725.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2121.6875
$ws.Range("I3").Value = 994.8
$ws.Range("J3").Value = 3999.8333
$ws.Range("K3").Value = 994.8
$ws.Range("L3").Value = 3999.8333
$ws.Range("M3").Value = -880.8
$ws.Range("N3").Value = -4227.8333
# Row 75
$ws.Range("H75").Value = 21824.75
$ws.Range("J75").Value = 4236
$ws.Range("L75").Value = 4236
$ws.Range("N75").Value = -6108
# Row 78
$ws.Range("H78").Value = 21824.75
$ws.Range("J78").Value = 4236
$ws.Range("L78").Value = 12708
$ws.Range("N78").Value = -22068

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2807.5
$ws.Range("J16").Value = 1677.5
$ws.Range("L16").Value = 1677.5
$ws.Range("N16").Value = -2251.5
# Row 54
$ws.Range("H54").Value = 37416.332
$ws.Range("I54").Value = 37416.332
$ws.Range("K54").Value = 37416.332
$ws.Range("M54").Value = -36758.332
# Row 58
$ws.Range("H58").Value = 896
$ws.Range("I58").Value = 896
$ws.Range("K58").Value = 896
$ws.Range("M58").Value = -693
# Row 113
$ws.Range("H113").Value = 2807.5
$ws.Range("J113").Value = 1677.5
$ws.Range("L113").Value = 1677.5
$ws.Range("N113").Value = -6017.5
# Row 136
$ws.Range("H136").Value = 896
$ws.Range("I136").Value = 896
$ws.Range("K136").Value = 2688
$ws.Range("M136").Value = -138

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
# Row 58
$ws.Range("H58").Value = 2466.3333
# Row 64
$ws.Range("H64").Value = 19500
$ws.Range("J64").Value = 19500
$ws.Range("L64").Value = 58500
$ws.Range("N64").Value = -59040
# Row 67
$ws.Range("H67").Value = 19500
$ws.Range("J67").Value = 19500
$ws.Range("L67").Value = 58500
$ws.Range("N67").Value = -60372
# Row 137
$ws.Range("H137").Value = 5146.143
$ws.Range("J137").Value = 5232.1665
$ws.Range("L137").Value = 15696.4995
$ws.Range("N137").Value = -25896.4995

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 162502
$ws.Range("I5").Value = 162502
$ws.Range("K5").Value = 162502
$ws.Range("M5").Value = -162390
# Row 113
$ws.Range("H113").Value = 2321.7778
$ws.Range("I113").Value = 1985.1428
$ws.Range("K113").Value = 1985.1428
$ws.Range("M113").Value = 184.8571999999999
# Row 132
$ws.Range("H132").Value = 2320.8572
$ws.Range("I132").Value = 2226.818
$ws.Range("J132").Value = 2665.6667
$ws.Range("K132").Value = 6680.454000000001
$ws.Range("L132").Value = 7997.000100000001
$ws.Range("M132").Value = -4150.454000000001
$ws.Range("N132").Value = -13057.0001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1945.4667
$ws.Range("I22").Value = 1804.375
$ws.Range("J22").Value = 2106.7144
$ws.Range("K22").Value = 1804.375
$ws.Range("L22").Value = 2106.7144
$ws.Range("M22").Value = -1509.375
$ws.Range("N22").Value = -2696.7144
# Row 27
$ws.Range("H27").Value = 1945.4667
$ws.Range("I27").Value = 1804.375
$ws.Range("J27").Value = 2106.7144
$ws.Range("K27").Value = 1804.375
$ws.Range("L27").Value = 2106.7144
$ws.Range("M27").Value = -1697.375
$ws.Range("N27").Value = -2320.7144
# Row 46
$ws.Range("H46").Value = 39615.11
$ws.Range("I46").Value = 85083.836
$ws.Range("K46").Value = 85083.836
$ws.Range("M46").Value = -84895.836
# Row 61
$ws.Range("H61").Value = 5462.636
$ws.Range("J61").Value = 7588
$ws.Range("L61").Value = 7588
$ws.Range("N61").Value = -7992
# Row 68
$ws.Range("H68").Value = 2149.25
$ws.Range("I68").Value = 1519.3
$ws.Range("J68").Value = 3199.1667
$ws.Range("K68").Value = 1519.3
$ws.Range("L68").Value = 3199.1667
$ws.Range("M68").Value = -770.3
$ws.Range("N68").Value = -4697.1667
# Row 71
$ws.Range("H71").Value = 2149.25
$ws.Range("I71").Value = 1519.3
$ws.Range("J71").Value = 3199.1667
$ws.Range("K71").Value = 7596.5
$ws.Range("L71").Value = 15995.8335
$ws.Range("M71").Value = -3852.5
$ws.Range("N71").Value = -23483.8335
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
# Row 100
$ws.Range("H100").Value = 2536.75
$ws.Range("I100").Value = 2152.7144
$ws.Range("J100").Value = 5225
$ws.Range("K100").Value = 2152.7144
$ws.Range("L100").Value = 5225
$ws.Range("M100").Value = -1611.7144
$ws.Range("N100").Value = -6307
# Row 113
$ws.Range("H113").Value = 5462.636
$ws.Range("J113").Value = 7588
$ws.Range("L113").Value = 7588
$ws.Range("N113").Value = -11928
# Row 127
$ws.Range("H127").Value = 33333
$ws.Range("J127").Value = 33333
$ws.Range("L127").Value = 33333
$ws.Range("N127").Value = -43253

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 20937.223
$ws.Range("I54").Value = 20937.223
$ws.Range("K54").Value = 20937.223
$ws.Range("M54").Value = -20417.223
# Row 75
$ws.Range("H75").Value = 88333.336
$ws.Range("I75").Value = 88333.336
$ws.Range("K75").Value = 88333.336
$ws.Range("M75").Value = -87397.336
# Row 78
$ws.Range("H78").Value = 88333.336
$ws.Range("I78").Value = 88333.336
$ws.Range("K78").Value = 265000.008
$ws.Range("M78").Value = -260320.008
# Row 126
$ws.Range("H126").Value = 4909.421
$ws.Range("I126").Value = 2981.182
$ws.Range("K126").Value = 8943.545999999998
$ws.Range("M126").Value = -6473.545999999998
